$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rescatables")

$ws.Cells.Item(2, 1).Value = 18330051920237
$ws.Cells.Item(2, 2).Value = "CANCINO"
$ws.Cells.Item(2, 3).Value = "GUERRA"
$ws.Cells.Item(2, 4).Value = "DANIEL"
$ws.Cells.Item(2, 5).Value = "TEMAS DE FILOSOFÍA"
$ws.Cells.Item(2, 6).Value = "6ALCV"
$ws.Cells.Item(2, 7).Value = 2
